{"js": "// Merge the two duplicate \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \" paragraphs\n// into a single paragraph whose run now reads \"\u0414\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u044d\u043b\u0435\u043c\u0435\u043d\u0442\" repeated\n// five times, and mark the run's language as Russian (ru-RU), matching\n// the target OOXML.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst first = paragraphs.items[0];\nconst second = paragraphs.items[1];\n\nconst newText = \"\u0414\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u044d\u043b\u0435\u043c\u0435\u043d\u0442\".repeat(5);\n\n// Replace the text of the first paragraph's run while keeping its\n// existing character formatting (rFonts/sz/szCs/color/kern untouched).\nfirst.insertText(newText, Word.InsertLocation.replace);\n\n// Tag the run with the Russian language id -> <w:lang w:val=\"ru-RU\"/>.\nconst firstRange = first.getRange();\nfirstRange.languageId = \"ru-RU\";\n\n// Drop the now-redundant second paragraph entirely.\nsecond.delete();\n\nawait context.sync();\n", "ps1": "# Merge the two duplicate \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \" paragraphs\n# into a single paragraph whose run reads \"\u0414\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u044d\u043b\u0435\u043c\u0435\u043d\u0442\" repeated\n# five times, tagged as Russian (ru-RU), matching the target OOXML.\n$d = $word.ActiveDocument\n\n$newText = \"\"\nfor ($i = 0; $i -lt 5; $i++) {\n    $newText = $newText + \"\u0414\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u044d\u043b\u0435\u043c\u0435\u043d\u0442\"\n}\n\n$firstPara = $d.Paragraphs.Item(1)\n$firstRange = $firstPara.Range\n$firstRange.Text = $newText\n$firstRange.LanguageID = \"ru-RU\"\n\n$secondPara = $d.Paragraphs.Item(2)\n$secondPara.Range.Delete()\n"}
